$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1582
$ws.Range("G2").Value = 208
$ws.Range("F3").Value = 968
$ws.Range("F4").Value = 945
$ws.Range("F5").Value = 542
$ws.Range("F6").Value = 8220
$ws.Range("F7").Value = 145
$ws.Range("F10").Value = 5902
$ws.Range("F11").Value = 590
$ws.Range("F14").Value = 8336
$ws.Range("F15").Value = 9791
$ws.Range("F17").Value = 977
$ws.Range("F18").Value = 4651
$ws.Range("F20").Value = 309
$ws.Range("F24").Value = 1242
$ws.Range("F25").Value = 151
$ws.Range("F26").Value = 1767
$ws.Range("F27").Value = 787
$ws.Range("F28").Value = 1048
$ws.Range("F29").Value = 423
$ws.Range("F30").Value = 1946
$ws.Range("F31").Value = 365
$ws.Range("F32").Value = 515
$ws.Range("F33").Value = 2436
$ws.Range("F35").Value = 129
$ws.Range("F36").Value = 1550
$ws.Range("F38").Value = 1331
$ws.Range("F39").Value = 22
$ws.Range("F40").Value = 833
$ws.Range("F42").Value = 212
$ws.Range("F45").Value = 537
$ws.Range("F49").Value = 4137

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F21").Value = 3

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 5514

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1582
$ws.Range("G3").Value = 208
$ws.Range("F4").Value = 968
$ws.Range("F5").Value = 945
$ws.Range("F6").Value = 542
$ws.Range("F7").Value = 8220
$ws.Range("F8").Value = 145
$ws.Range("F12").Value = 5902
$ws.Range("F13").Value = 590
$ws.Range("F14").Value = 8336
$ws.Range("F15").Value = 9791
$ws.Range("F18").Value = 977
$ws.Range("F19").Value = 4651
$ws.Range("F21").Value = 309
$ws.Range("F25").Value = 1242
$ws.Range("F26").Value = 151
$ws.Range("F27").Value = 1767
$ws.Range("F28").Value = 787
$ws.Range("F29").Value = 1048
$ws.Range("F30").Value = 423
$ws.Range("F32").Value = 1946
$ws.Range("F33").Value = 365
$ws.Range("F34").Value = 2436
$ws.Range("F38").Value = 833
$ws.Range("F42").Value = 212
$ws.Range("F45").Value = 537
$ws.Range("F48").Value = 4137
